# Regenerate save_data: update column G ("K") values for rows 2-22
# to reflect new computed K counts (replacing old Strike#-based values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 2
    4  = 4
    5  = 4
    6  = 5
    7  = 3
    8  = 3
    9  = 7
    10 = 4
    11 = 5
    12 = 4
    13 = 4
    14 = 8
    15 = 5
    16 = 6
    17 = 8
    18 = 9
    19 = 1
    20 = 2
    21 = 1
    22 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
